$wb = $excel.ActiveWorkbook

# Insert a new "is_targeted list" sheet right after "analyte_class list"
# (shifts gdna_fragmenta...assurance list / dna_assay_input_unit list /
#  library_layout list / library_final_yield_unit list down by one tab).
$afterSheet = $wb.Worksheets.Item("analyte_class list")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "is_targeted list"

# Populate it with the TRUE / FALSE list values (force text, not boolean).
$newSheet.Range("A1").Value = "'TRUE"
$newSheet.Range("A2").Value = "'FALSE"

# Point column N's data validation at the new list sheet instead of the
# inline "TRUE,FALSE" formula, and update the error text to match.
$ws1 = $wb.Worksheets.Item("Export as TSV")
$nRange = $ws1.Range("N2:N1048576")
$nRange.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$nRange.Validation.ErrorTitle = "Value must come from list"
$nRange.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

Write-Output "ok"
